$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp string in A1 (22:45 -> 23:15)
$ws.Range("A1").Value = "Datos actualizados a 19 de Marzo de 2020 a las 23:15"

# Swap the province names in rows 52 and 53 (Melilla <-> Huelva), keeping
# their numeric data (Casos totales/activos/Recuperados/Muertes) untouched.
$ws.Range("A52").Value = "Huelva"
$ws.Range("A53").Value = "Melilla"

# Update Asturias (row 14) Recuperados/Muertes figures.
$ws.Range("D14").Value = 285
$ws.Range("E14").Value = 3
